$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6435753333333333
$ws.Range("H2").Value = 1.930726
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.7517265
$ws.Range("N2").Value = 1.503453
$ws.Range("O2").Value = 0.3430414969595709
$ws.Range("P2").Value = 0.2862361877440628
$ws.Range("Q2").Value = 0.4837926328129999
$ws.Range("R2").Value = 2.902755796878
$ws.Range("S2").Value = 0.3430414969595709
$ws.Range("T2").Value = 0.2862361877440628

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6435753333333333
$ws.Range("H3").Value = 1.930726
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.2592623333333333
$ws.Range("N3").Value = 0.777787
$ws.Range("O3").Value = 0.1183112992982127
$ws.Range("P3").Value = 0.1480796444962971
$ws.Range("Q3").Value = 0.1668548425957778
$ws.Range("R3").Value = 1.501693583362
$ws.Range("S3").Value = 0.1183112992982127
$ws.Range("T3").Value = 0.1480796444962971

# Row 4 (Target cluster: M1)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6435753333333333
$ws.Range("H4").Value = 1.930726
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1423686666666667
$ws.Range("N4").Value = 0.427106
$ws.Range("O4").Value = 0.06496825711674591
$ws.Range("P4").Value = 0.08131494180570706
$ws.Range("Q4").Value = 0.09162496210622222
$ws.Range("R4").Value = 0.8246246589559999
$ws.Range("S4").Value = 0.06496825711674591
$ws.Range("T4").Value = 0.08131494180570706

# Row 5 (Target cluster: M2)
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.6435753333333333
$ws.Range("H5").Value = 1.930726
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.171105
$ws.Range("N5").Value = 0.513315
$ws.Range("O5").Value = 0.07808174294409917
$ws.Range("P5").Value = 0.09772791614493011
$ws.Range("Q5").Value = 0.11011895741
$ws.Range("R5").Value = 0.9910706166899999
$ws.Range("S5").Value = 0.07808174294409917
$ws.Range("T5").Value = 0.09772791614493011

# Row 6 (Target cluster: Neutro)
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.6435753333333333
$ws.Range("H6").Value = 1.930726
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.2970403333333334
$ws.Range("N6").Value = 0.8911210000000001
$ws.Range("O6").Value = 0.1355508427653363
$ws.Range("P6").Value = 0.1696568352044773
$ws.Range("Q6").Value = 0.1911678315384444
$ws.Range("R6").Value = 1.720510483846
$ws.Range("S6").Value = 0.1355508427653363
$ws.Range("T6").Value = 0.1696568352044773

# Row 7 (Target cluster: sCs)
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.6435753333333333
$ws.Range("H7").Value = 1.930726
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.5698545
$ws.Range("N7").Value = 1.139709
$ws.Range("O7").Value = 0.260046360916035
$ws.Range("P7").Value = 0.2169844746045258
$ws.Range("Q7").Value = 0.366744299789
$ws.Range("R7").Value = 2.200465798734
$ws.Range("S7").Value = 0.260046360916035
$ws.Range("T7").Value = 0.2169844746045258
